# Edit script: update model predictions (IPC PO), DELTA and DELTA^2 columns
# so that the predicted value (C) is now constant 0, DELTA (D) becomes the
# negation of the observed value (B), and DELTA^2 (E) becomes B^2.
# Also refresh the TOTAL (row 52) and MSE (row 53) summary rows accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$totalDelta = 0
$totalDeltaSq = 0
$count = 0

for ($r = 2; $r -le 51; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2

    $c = 0
    $d = 0 - $b
    $e = $d * $d

    $ws.Cells.Item($r, 3).Value2 = $c
    $ws.Cells.Item($r, 4).Value2 = $d
    $ws.Cells.Item($r, 5).Value2 = $e

    $totalDelta = $totalDelta + $d
    $totalDeltaSq = $totalDeltaSq + $e
    $count = $count + 1
}

# Row 52: TOTAL
$ws.Cells.Item(52, 3).Value2 = $totalDelta
$ws.Cells.Item(52, 5).Value2 = $totalDeltaSq

# Row 53: MSE
$ws.Cells.Item(53, 5).Value2 = $totalDeltaSq / $count

Write-Host "Updated rows 2-53 with new predictions (constant 0), deltas and squared deltas."
